$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "NASA" Young's Modulus table (columns H:L, rows 1-7)
# ---------------------------------------------------------------------------

$ws.Range("H1").Value = "NASA"
$ws.Range("I1").Value = "Temp (K)"
$ws.Range("J1").Value = "Young's Modulus"

$ws.Range("I2").Value = 77.594444444444434
$ws.Range("J2").Value = 79461934113.818359

$ws.Range("I3").Value = 294.26111111111106
$ws.Range("J3").Value = 71877714805.774963

$ws.Range("I4").Value = 373.15
$ws.Range("J4").Value = 71188240323.225555

$ws.Range("I5").Value = 473.15
$ws.Range("J5").Value = 58950068257.973709

$ws.Range("I6").Value = 573.15
$ws.Range("J6").Value = 37748727919.579651

$ws.Range("I7").Value = 673.15
$ws.Range("J7").Value = 19535110338.89967

Write-Host "values set"

# ---------------------------------------------------------------------------
# Borders: draw a box outline around H1:J7 (left on H, top on row1, right on
# J1 only, bottom on row7) to match the author's header/table styling.
# ---------------------------------------------------------------------------

# Left edge of the box (column H, rows 1-7)
$ws.Range("H1:H7").Borders.Item(7).LineStyle = 1
$ws.Range("H1:H7").Borders.Item(7).Weight = 2

# Top edge of the box (row 1, columns H:J)
$ws.Range("H1:J1").Borders.Item(8).LineStyle = 1
$ws.Range("H1:J1").Borders.Item(8).Weight = 2

# Right edge of the box (top-right corner cell only, J1)
$ws.Range("J1").Borders.Item(10).LineStyle = 1
$ws.Range("J1").Borders.Item(10).Weight = 2

# Bottom edge of the box (row 7, columns H:J)
$ws.Range("H7:J7").Borders.Item(9).LineStyle = 1
$ws.Range("H7:J7").Borders.Item(9).Weight = 2

Write-Host "borders set"
